# Update cryptos list values (Price and Volume(1h) columns) per scraped data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.072.62'
$ws.Range("E2").Value = '  +2.24%  '
$ws.Range("D3").Value = '3.819.48'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = "'628.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.30%  '
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("D7").Value = '3.816.87'
$ws.Range("E7").Value = '  +0.97%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("E9").Value = '  +1.18%  '
$ws.Range("D10").Value = "'0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.54%  '
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("D12").Value = "'6.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.20%  '
$ws.Range("D13").Value = "'0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.31%  '
$ws.Range("D14").Value = "'36.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("D15").Value = '4.457.45'
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("D16").Value = '3.932.58'
$ws.Range("E16").Value = '  +2.98%  '
$ws.Range("D17").Value = '69.070.66'
$ws.Range("E17").Value = '  +2.12%  '
$ws.Range("E18").Value = '  -1.26%  '
$ws.Range("D19").Value = "'7.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.80%  '
$ws.Range("D21").Value = "'465.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.07%  '
$ws.Range("D22").Value = "'9.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.40%  '
$ws.Range("D23").Value = "'0.708"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.98%  '
$ws.Range("E24").Value = '  +4.96%  '
$ws.Range("E25").Value = '  +1.68%  '
$ws.Range("D26").Value = "'11.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  +3.27%  '
$ws.Range("E28").Value = '  +0.55%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").Value = '3.967.57'
$ws.Range("E30").Value = '  +0.89%  '
$ws.Range("E31").Value = '  +1.73%  '
$ws.Range("E32").Value = '  +1.82%  '
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("D34").Value = "'29.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.11%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  +1.48%  '
$ws.Range("D37").Value = "'0.102"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.14%  '
$ws.Range("D38").Value = "'0.149"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.70%  '
$ws.Range("D39").Value = "'3.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.97%  '
$ws.Range("E40").Value = '  +3.39%  '
$ws.Range("E41").Value = '  -0.65%  '
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = "'157.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.59%  '
$ws.Range("E45").Value = '  +5.85%  '
$ws.Range("E46").Value = '  +1.14%  '
$ws.Range("D47").Value = "'46.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.56%  '
$ws.Range("D48").Value = "'42.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.83%  '
$ws.Range("E49").Value = '  +1.75%  '
$ws.Range("E50").Value = '  +3.15%  '
$ws.Range("D51").Value = "'0.000280"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +13.97%  '
